$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.649.86"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "2.528.80"
$ws.Range("E3").Value = "  -1.85%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("E7").Value = "  -1.19%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -2.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.87%  "

$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("D14").Value = "2.916.77"
$ws.Range("E14").Value = "  -1.88%  "

$ws.Range("E15").Value = "  -2.78%  "

$ws.Range("D16").Value = "2.509.57"
$ws.Range("E16").Value = "  -4.14%  "

$ws.Range("E17").Value = "  -4.39%  "

$ws.Range("D18").Value = "42.644.59"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.40%  "

$ws.Range("D20").Value = "0.0₃0950"
$ws.Range("E20").Value = "  -1.86%  "

$ws.Range("E21").Value = "  -3.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.26"
$ws.Range("D23").Style = "Normal"

$ws.Range("E24").Value = "  -3.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.21%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.14%  "

$ws.Range("E28").Value = "  -2.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.47%  "

$ws.Range("E33").Value = "  +10.23%  "

$ws.Range("E34").Value = "  -1.48%  "

$ws.Range("E35").Value = "  -2.87%  "

$ws.Range("E36").Value = "  -8.98%  "

$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.59%  "

$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.25%  "

$ws.Range("E39").Value = "  -1.50%  "

$ws.Range("E40").Value = "  -0.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.82%  "

$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("E44").Value = "  -1.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.94%  "

$ws.Range("D46").Value = "2.007.76"
$ws.Range("E46").Value = "  +0.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("D48").Value = "2.770.83"
$ws.Range("E48").Value = "  -1.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.190"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.08%  "
